# Auto-generated Excel COM-interop script
# Applies the cell-content updates for the "cryptos" worksheet
# (price / volume refresh + two rank swaps: Hedera<->LidoDAOToken, Aave->TrustWalletToken)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.961.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'2.217.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'291.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'87.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.467"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").Value = "'30.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "'50.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.26%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'6.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "'2.560.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "'13.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "'2.245.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "'0.731"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'39.913.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "'5.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'65.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'237.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'2.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "'23.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'9.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -6.68%  "
$ws.Range("D31").Value = "'156.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("D32").Value = "'31.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "'4.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0712"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'0.0990"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "'15.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("D42").Value = "'2.105.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").Value = "'3.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "'17.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Value = "'9.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").Value = "'1.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.28%  "
$ws.Range("D48").Value = "'2.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("D49").Value = "'2.431.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.55%  "
